# Rename the sheet to match the new test case name.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "AddCustomerTest"

# Header row - written in the same "column by column" order the
# original data-driven test fixture used (lastName/firstName/postCode
# first, "alertText" appended afterwards for the CustomListeners work).
$ws.Range("B1").Value = "lastName"
$ws.Range("A1").Value = "firstName"
$ws.Range("C1").Value = "postCode"

# Data row for the AddCustomerTest scenario.
$ws.Range("A2").Value = "Sohaib"
$ws.Range("B2").Value = "Majeed"
$ws.Range("C2").Value = "123wp"

# New "alertText" column added for the ReportNG / custom listener work.
$ws.Range("D1").Value = "alertText"
$ws.Range("D2").Value = "Customer added successfully"

# Header row rendered in bold.
$ws.Range("A1:D1").Font.Bold = $true

# Size the columns to fit their (now wider) contents (bestFit/autofit).
$ws.Columns.Item(1).ColumnWidth = 9
$ws.Columns.Item(2).ColumnWidth = 8.666666666666666
$ws.Columns.Item(3).ColumnWidth = 8.666666666666666
$ws.Columns.Item(4).ColumnWidth = 42.166666666666664

# Leave the cursor where data entry finished.
[void]$ws.Range("D3").Select()

# Portrait page setup for printing the sheet.
$ws.PageSetup.Orientation = 1
